$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Fix two author names in row 2 (A2): add DOM/Banner author ids
#    "A.O. King" -> "Amanda R. King"
#    "Suraj Venna" -> "Suraj S. Venna"
$a2 = $ws.Range("A2").Value2
$a2 = $a2.Replace("A.O. King", "Amanda R. King")
$a2 = $a2.Replace("Suraj Venna", "Suraj S. Venna")
$ws.Range("A2").Value = "'" + $a2
$ws.Range("A2").Style = "Normal"

# 2) Swap the data of row 3 and row 4 (all columns A:Q) - a new
#    article row was inserted ahead of the previously-first article,
#    which (since rows were appended rather than truly inserted)
#    shows up as the two rows' contents being exchanged.
$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q")

$row3vals = @{}
$row4vals = @{}
foreach ($col in $cols) {
    $row3vals[$col] = $ws.Range($col + "3").Value2
    $row4vals[$col] = $ws.Range($col + "4").Value2
}

foreach ($col in $cols) {
    $ws.Range($col + "3").Value = "'" + $row4vals[$col]
    $ws.Range($col + "4").Value = "'" + $row3vals[$col]
}

$ws.Range("A3:Q4").Style = "Normal"
